$d = $word.ActiveDocument

# Update the date/day heading at the top of the document
$d.Content.Find.Execute("2024-01-25 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-01-26 Friday", 2)

# Update the division problems in the table. Cells are addressed directly by
# (row, column) position so that duplicate/overlapping old & new values
# (e.g. "24÷5=" and "15÷5=" appear both as a source and a result) are each
# replaced exactly once, without cascading into each other.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "44÷4="
$t.Cell(1, 2).Range.Text  = "15÷5="
$t.Cell(1, 3).Range.Text  = "16÷4="
$t.Cell(1, 4).Range.Text  = "24÷7="
$t.Cell(1, 5).Range.Text  = "70÷3="

$t.Cell(5, 1).Range.Text  = "24÷5="
$t.Cell(5, 2).Range.Text  = "55÷8="
$t.Cell(5, 3).Range.Text  = "89÷8="
$t.Cell(5, 4).Range.Text  = "68÷2="
$t.Cell(5, 5).Range.Text  = "40÷2="

$t.Cell(9, 1).Range.Text  = "57÷5="
$t.Cell(9, 2).Range.Text  = "79÷2="
$t.Cell(9, 3).Range.Text  = "77÷5="
$t.Cell(9, 4).Range.Text  = "21÷3="
$t.Cell(9, 5).Range.Text  = "53÷5="

$t.Cell(13, 1).Range.Text = "10÷8="
$t.Cell(13, 2).Range.Text = "19÷4="
$t.Cell(13, 3).Range.Text = "31÷2="
$t.Cell(13, 4).Range.Text = "10÷7="
$t.Cell(13, 5).Range.Text = "64÷9="

$t.Cell(17, 1).Range.Text = "68÷8="
$t.Cell(17, 2).Range.Text = "25÷7="
$t.Cell(17, 3).Range.Text = "85÷7="
$t.Cell(17, 4).Range.Text = "56÷4="
$t.Cell(17, 5).Range.Text = "30÷5="
